# Add a "Supervisors" sheet (with per-supervisor max-projects / max-students
# caps) ahead of "Projects", and record each project's supervisor on the
# "Projects" sheet. Data content of Student_preferences / Supervisor_preferences
# is untouched.

$wb = $excel.ActiveWorkbook

# --- New "Supervisors" sheet, inserted before "Projects" --------------------
$projectsBefore = $wb.Worksheets.Item("Projects")
$supervisors = $wb.Worksheets.Add($projectsBefore)
$supervisors.Name = "Supervisors"

# NOTE: inserting a sheet shifts collection indices, so re-fetch "Projects"
# by name afterwards rather than reusing the handle obtained beforehand.
$projects = $wb.Worksheets.Item("Projects")

$supervisors.Range("A1").Value = "Supervisor"
$supervisors.Range("B1").Value = "Max_number_of_projects"
$supervisors.Range("C1").Value = "Max_number_of_students"
$supervisors.Range("A2").Value = "Dr Smith"
$supervisors.Columns.Item(2).ColumnWidth = 23.16

# --- "Projects" sheet: rename header, add Supervisor column -----------------
$projects.Range("A1").Value = "Project"
$projects.Range("C1").Value = "Supervisor"
for ($r = 2; $r -le 11; $r++) {
    $projects.Cells.Item($r, 3).Value = "Dr Smith"
}
$projects.Columns.Item(2).ColumnWidth = 24.07

# --- Selection / active-sheet bookkeeping ------------------------------------
$supervisors.Activate()
[void]$supervisors.Range("B2").Select()

$projects.Activate()
[void]$projects.Range("C2:C11").Select()
